# Updates cryptos list values (price & 1h volume change) per scheduled data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.059.84"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.423.93"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'406.94"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'131.77"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").Value = "'0.594"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.691"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("E10").Value = "  +6.82%  "
$ws.Range("D11").Value = "'41.89"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D13").Value = "'19.81"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'8.37"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "3.457.45"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "'11.62"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "62.043.31"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  +9.33%  "
$ws.Range("D20").Value = "'3.17"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").Value = "'83.81"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "'311.56"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'12.78"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").Value = "'4.70"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = "  +4.12%  "
$ws.Range("E29").Value = "  +5.48%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'43.64"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'0.0485"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'51.78"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'2.99"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.318"
$ws.Range("E39").Value = "  +11.91%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'3.31"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "'143.81"
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").Value = "'16.84"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'21.17"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").Value = "2.102.26"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("E51").Value = "  +17.36%  "
